$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.323.81'
$ws.Range('E2').Value = '  -3.39%  '
$ws.Range('D3').Value = '1.934.17'
$ws.Range('E3').Value = '  -3.70%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = "'249.33"
$ws.Range('E5').Value = '  -3.83%  '
$ws.Range('D6').Value = "'0.7236"
$ws.Range('E6').Value = '  -8.25%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').Value = "'0.3288"
$ws.Range('E8').Value = '  -8.69%  '
$ws.Range('D9').Value = "'27.84"
$ws.Range('E9').Value = '  -3.06%  '
$ws.Range('D10').Value = "'0.06863"
$ws.Range('E10').Value = '  -3.05%  '
$ws.Range('D11').Value = "'0.8062"
$ws.Range('E11').Value = '  -5.50%  '
$ws.Range('D12').Value = "'0.08080"
$ws.Range('E12').Value = '  -0.19%  '
$ws.Range('D13').Value = '1.929.94'
$ws.Range('E13').Value = '  -3.93%  '
$ws.Range('D14').Value = "'5.415"
$ws.Range('E14').Value = '  -3.77%  '
$ws.Range('D15').Value = "'94.89"
$ws.Range('E15').Value = '  -6.53%  '
$ws.Range('D16').Value = "'14.54"
$ws.Range('E16').Value = '  -1.86%  '
$ws.Range('D17').Value = '30.320.62'
$ws.Range('E17').Value = '  -3.41%  '
$ws.Range('D18').Value = "'252.72"
$ws.Range('E18').Value = '  -8.48%  '
$ws.Range('D19').Value = "'0.000008116"
$ws.Range('E19').Value = '  +2.27%  '
$ws.Range('D20').Value = "'5.830"
$ws.Range('E20').Value = '  -2.07%  '
$ws.Range('D21').Value = '2.185.75'
$ws.Range('E21').Value = '  -3.39%  '
$ws.Range('E22').Value = '  +0.09%  '
$ws.Range('E23').Value = '  +0.00%  '
$ws.Range('D24').Value = "'6.874"
$ws.Range('E24').Value = '  -4.71%  '
$ws.Range('D25').Value = "'9.726"
$ws.Range('E25').Value = '  -3.90%  '
$ws.Range('D26').Value = "'159.55"
$ws.Range('E26').Value = '  -3.03%  '
$ws.Range('D27').Value = "'2.402"
$ws.Range('E27').Value = '  +0.46%  '
$ws.Range('E28').Value = '  -4.77%  '
$ws.Range('D29').Value = "'0.1340"
$ws.Range('E29').Value = '  -11.21%  '
$ws.Range('D30').Value = "'1.560"
$ws.Range('E30').Value = '  -4.05%  '
$ws.Range('D31').Value = "'1.339"
$ws.Range('E31').Value = '  -1.43%  '
$ws.Range('E32').Value = '  -5.14%  '
$ws.Range('D33').Value = "'4.191"
$ws.Range('E33').Value = '  -4.97%  '
$ws.Range('D34').Value = "'0.05109"
$ws.Range('E34').Value = '  -2.59%  '
$ws.Range('D35').Value = "'1.223"
$ws.Range('E35').Value = '  -0.39%  '
$ws.Range('D36').Value = "'0.7418"
$ws.Range('E36').Value = '  -3.32%  '
$ws.Range('D37').Value = "'2.750"
$ws.Range('E37').Value = '  -2.20%  '
$ws.Range('E38').Value = '  -2.31%  '
$ws.Range('E39').Value = '  -4.10%  '
$ws.Range('D40').Value = "'6.619"
$ws.Range('E40').Value = '  -1.41%  '
$ws.Range('D41').Value = "'79.02"
$ws.Range('E41').Value = '  -2.37%  '
$ws.Range('D42').Value = "'0.4469"
$ws.Range('E42').Value = '  -6.09%  '
$ws.Range('D43').Value = "'1.998"
$ws.Range('E43').Value = '  -8.81%  '
$ws.Range('E44').Value = '  +0.08%  '
$ws.Range('D45').Value = "'0.8358"
$ws.Range('E45').Value = '  -2.74%  '
$ws.Range('D46').Value = "'102.13"
$ws.Range('E46').Value = '  -2.88%  '
$ws.Range('D47').Value = "'9.779"
$ws.Range('E47').Value = '  -2.00%  '
$ws.Range('D48').Value = "'7.321"
$ws.Range('E48').Value = '  -4.98%  '
$ws.Range('D49').Value = "'36.65"
$ws.Range('E49').Value = '  -0.79%  '
$ws.Range('D50').Value = "'0.05957"
$ws.Range('E50').Value = '  -0.48%  '
$ws.Range('E51').Value = '  -0.42%  '
